$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Execute code from Eclipse or JAR file" -> append a new run with
#    " (JAR file found in the “program” directory)"
#    A temporary bookmark is used purely to stop the host from
#    re-merging the freshly-typed text back into the preceding run
#    (it is deleted again immediately, leaving two plain <w:r> runs,
#    exactly like the target markup).
# ---------------------------------------------------------------------
$rFind = $d.Content
$rFind.Find.Execute("Execute code from Eclipse or JAR file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rFind.Collapse(0)
$d.Bookmarks.Add("ZZZtempsplit", $rFind) | Out-Null
$rFind.InsertAfter(" (JAR file found in the “program” directory)")
$d.Bookmarks("ZZZtempsplit").Delete()

# ---------------------------------------------------------------------
# 2) Tag the run that hosts the screenshot drawing as <w:noProof/>.
# ---------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$shp.Range.NoProofing = 1

# ---------------------------------------------------------------------
# 3) Drop the "TP-2/" prefix in front of "Student Input File.xls".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("TP-2/Student Input File.xls", $true, $false, $false, $false, $false, $true, 1, $false, "Student Input File.xls", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Drop the "TP-2/" prefix in front of "Project Input File.xls" and
#    split that sentence into two runs with an (empty) "_GoBack"
#    bookmark sitting between them, moving the bookmark away from its
#    old home at the end of the "build path" bullet.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Content.Find.Execute("TP-2/Project Input File.xls", $true, $false, $false, $false, $false, $true, 1, $false, "Project Input File.xls", 2) | Out-Null

$rSplit = $d.Content
$rSplit.Find.Execute("Press the “Read Projects” button to read a VALID projects file.  The provided “", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rSplit.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rSplit) | Out-Null
